# MHTG_RNAV_FINALES.xlsx - minor arc center correction
# - Rename "NORTH RNAV (RNP) RWY 02" -> "NORTH1B"
# - Rename "SOUTH RNAV (RNP) RWY 02" -> "SOUTH1B"
# - Make the NORTH1B sheet the active tab/view (zoom 80%, selection Y7)
# - Widen a few columns on the NORTH1B sheet
# - Previously selected sheet (MHTG 2-39.10 RNAV (RNP) RW20) loses the
#   "tabSelected" flag automatically once a different sheet is activated.

$wb = $excel.ActiveWorkbook

# --- Rename worksheets (by original name, order-independent) ---
foreach ($sheet in $wb.Worksheets) {
    if ($sheet.Name -eq "NORTH RNAV (RNP) RWY 02") {
        $sheet.Name = "NORTH1B"
    }
    elseif ($sheet.Name -eq "SOUTH RNAV (RNP) RWY 02") {
        $sheet.Name = "SOUTH1B"
    }
}

$wsNorth = $wb.Worksheets.Item("NORTH1B")

# --- Make NORTH1B the active sheet/tab, with its own zoom + selection ---
$wsNorth.Activate()
$excel.ActiveWindow.Zoom = 80
$wsNorth.Range("Y7").Select()

# --- Column width tweaks on NORTH1B ---
$wsNorth.Columns.Item(11).ColumnWidth = 9.9166666666666667   # K -> width 10.75
$wsNorth.Columns.Item(15).ColumnWidth = 11.0833333333333333  # O -> width 12
$wsNorth.Columns.Item(16).ColumnWidth = 5.9166666666666667   # P -> width 6.75
